$wb = $excel.ActiveWorkbook

# This script applies the scheduled-runner price/profit refresh described in the
# commit: numeric recalculated values for currentAveragePrice / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ (and related) columns across several sheets, plus a couple of
# cells that appear/disappear because their row no longer has a computed profit figure.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 134.14285
$ws.Range("I33").Value = 159.33333
$ws.Range("K33").Value = 159.33333
$ws.Range("M33").Value = 69.66667000000001
$ws.Range("H39").Value = 1280.1428
$ws.Range("I39").Value = 503.5
$ws.Range("K39").Value = 1510.5
$ws.Range("M39").Value = -1214.5
$ws.Range("H40").Value = 1070.8918
$ws.Range("I40").Value = 1066.4667
$ws.Range("K40").Value = 1066.4667
$ws.Range("M40").Value = -891.4666999999999
$ws.Range("H112").Value = 1401.4706
$ws.Range("J112").Value = 1621.1538
$ws.Range("L112").Value = 4863.4614
$ws.Range("N112").Value = -7079.4614

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 1162.931
$ws.Range("I137").Value = 885.5217
$ws.Range("J137").Value = 2226.3333
$ws.Range("K137").Value = 2656.5651
$ws.Range("L137").Value = 6678.999899999999
$ws.Range("M137").Value = -106.5650999999998
$ws.Range("N137").Value = -11778.9999
$ws.Range("H138").Value = 1800.0426
$ws.Range("J138").Value = 3772.75
$ws.Range("L138").Value = 11318.25
$ws.Range("N138").Value = -21598.25
$ws.Range("H141").Value = 684044.3
$ws.Range("I141").Value = 737521.5
$ws.Range("J141").Value = 6666.3335
$ws.Range("K141").Value = 2212564.5
$ws.Range("L141").Value = 19999.0005
$ws.Range("M141").Value = -2207384.5
$ws.Range("N141").Value = -30359.0005
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("H32").Value = 3069.7356
$ws.Range("I32").Value = 2508.8513
$ws.Range("K32").Value = 2508.8513
$ws.Range("M32").Value = -2221.8513
$ws.Range("H39").Value = 3999.5
$ws.Range("I39").Value = 3999.5
$ws.Range("K39").Value = 3999.5
$ws.Range("M39").Value = -3479.5
$ws.Range("H61").Value = 2511.5908
$ws.Range("I61").Value = 1997.0555
$ws.Range("J61").Value = 4827
$ws.Range("K61").Value = 1997.0555
$ws.Range("L61").Value = 4827
$ws.Range("M61").Value = -1785.0555
$ws.Range("N61").Value = -5251
$ws.Range("H74").Value = 1225.9474
$ws.Range("I74").Value = 1020.6539
$ws.Range("K74").Value = 1020.6539
$ws.Range("M74").Value = -146.6539
$ws.Range("H77").Value = 1225.9474
$ws.Range("I77").Value = 1020.6539
$ws.Range("K77").Value = 5103.2695
$ws.Range("M77").Value = -735.2695000000003
$ws.Range("M26").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 1746
$ws.Range("I132").Value = 1266.8572
$ws.Range("J132").Value = 2704.2856
$ws.Range("K132").Value = 3800.5716
$ws.Range("L132").Value = 8112.8568
$ws.Range("M132").Value = -1270.5716
$ws.Range("N132").Value = -13172.8568
$ws.Range("H136").Value = 2511.5908
$ws.Range("I136").Value = 1997.0555
$ws.Range("J136").Value = 4827
$ws.Range("K136").Value = 5991.166499999999
$ws.Range("L136").Value = 14481
$ws.Range("M136").Value = -3441.166499999999
$ws.Range("N136").Value = -19581
$ws.Range("H68").Value = 29000
$ws.Range("H71").Value = 29000
$ws.Range("H82").Value = 14253.5
$ws.Range("I82").Value = 14253.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 14253.5
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 14253.5
$ws.Range("I85").Value = 14253.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 14253.5
$ws.Range("L85").Value = 0
$ws.Range("H105").Value = 2495.6086
$ws.Range("I105").Value = 2419.95
$ws.Range("K105").Value = 2419.95
$ws.Range("M105").Value = -672.9499999999998
$ws.Range("M82").Value = -13870.5
$ws.Range("M85").Value = -12927.5
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1400
$ws.Range("I16").Value = 800
$ws.Range("K16").Value = 800
$ws.Range("M16").Value = -513
$ws.Range("H31").Value = 1871.9546
$ws.Range("I31").Value = 1441.5
$ws.Range("J31").Value = 2117.9285
$ws.Range("K31").Value = 1441.5
$ws.Range("L31").Value = 2117.9285
$ws.Range("M31").Value = -1146.5
$ws.Range("N31").Value = -2707.9285
$ws.Range("H34").Value = 1871.9546
$ws.Range("I34").Value = 1441.5
$ws.Range("J34").Value = 2117.9285
$ws.Range("K34").Value = 1441.5
$ws.Range("L34").Value = 2117.9285
$ws.Range("M34").Value = -1239.5
$ws.Range("N34").Value = -2521.9285
$ws.Range("H52").Value = 64280
$ws.Range("J52").Value = 64280
$ws.Range("L52").Value = 64280
$ws.Range("N52").Value = -64868
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 800
$ws.Range("K113").Value = 800
$ws.Range("M113").Value = 1370

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2249.6667
$ws.Range("I48").Value = 1500
$ws.Range("K48").Value = 4500
$ws.Range("H116").Value = 71429840
$ws.Range("I116").Value = 774
$ws.Range("K116").Value = 2322
$ws.Range("M116").Value = 1120
$ws.Range("H117").Value = 815
$ws.Range("I117").Value = 494.25
$ws.Range("K117").Value = 1482.75
$ws.Range("M117").Value = 1959.25
$ws.Range("M48").Value = -4250

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 1299
$ws.Range("I130").Value = 1299
$ws.Range("K130").Value = 3897
$ws.Range("M130").Value = 1123
$ws.Range("H131").Value = 12336.159
$ws.Range("J131").Value = 14857.263
$ws.Range("L131").Value = 44571.789
$ws.Range("N131").Value = -54651.789
$ws.Range("H132").Value = 1070.0714
$ws.Range("I132").Value = 899.5
$ws.Range("K132").Value = 8095.5
$ws.Range("M132").Value = -5565.5
$ws.Range("H2").Value = 51.842106
$ws.Range("I2").Value = 11.5
$ws.Range("K2").Value = 11.5
$ws.Range("M2").Value = 101.5
$ws.Range("H102").Value = 2907.077
$ws.Range("I102").Value = 2981.182
$ws.Range("J102").Value = 2499.5
$ws.Range("K102").Value = 2981.182
$ws.Range("L102").Value = 2499.5
$ws.Range("M102").Value = -1359.182
$ws.Range("N102").Value = -5743.5
$ws.Range("H113").Value = 1587.1428
$ws.Range("I113").Value = 1752.5
$ws.Range("K113").Value = 1752.5
$ws.Range("M113").Value = 417.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 714300.3
$ws.Range("I132").Value = 1040926.94
$ws.Range("J132").Value = 3407.1765
$ws.Range("K132").Value = 3122780.82
$ws.Range("L132").Value = 10221.5295
$ws.Range("M132").Value = -3120250.82
$ws.Range("N132").Value = -15281.5295
$ws.Range("H134").Value = 24149.889
$ws.Range("J134").Value = 24149.889
$ws.Range("L134").Value = 72449.667
$ws.Range("N134").Value = -77519.667
$ws.Range("H22").Value = 4140.143
$ws.Range("I22").Value = 4000.3333
$ws.Range("J22").Value = 4245
$ws.Range("K22").Value = 4000.3333
$ws.Range("L22").Value = 4245
$ws.Range("M22").Value = -3705.3333
$ws.Range("N22").Value = -4835
$ws.Range("H27").Value = 4140.143
$ws.Range("I27").Value = 4000.3333
$ws.Range("J27").Value = 4245
$ws.Range("K27").Value = 4000.3333
$ws.Range("L27").Value = 4245
$ws.Range("M27").Value = -3893.3333
$ws.Range("N27").Value = -4459
$ws.Range("H61").Value = 1700.6428
$ws.Range("I61").Value = 1405.3636
$ws.Range("K61").Value = 1405.3636
$ws.Range("M61").Value = -1203.3636
$ws.Range("H100").Value = 2050
$ws.Range("I100").Value = 1650
$ws.Range("K100").Value = 1650
$ws.Range("M100").Value = -1109
$ws.Range("H113").Value = 1700.6428
$ws.Range("I113").Value = 1405.3636
$ws.Range("K113").Value = 1405.3636
$ws.Range("M113").Value = 764.6364000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 679
$ws.Range("I100").Value = 521.125
$ws.Range("K100").Value = 1042.25
$ws.Range("M100").Value = -501.25

Write-Output "Applied scheduled Sheets update."